# Automatische test-sync: 2025-07-27 19:45:50
# Append the new "Testmail #16" log row (row 18) to the "Logs" sheet,
# extend the conditional-formatting ranges to cover it, and bump the
# "Intern verzoek / Actie voor medewerker" counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- New row 18 on "Logs" ---------------------------------------------
$logs.Cells.Item(18, 1).Value = "Wil je dit even doorsturen?"
$logs.Cells.Item(18, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(18, 3).Value = "Testmail #16: Wil je dit even doorsturen?"
$logs.Cells.Item(18, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(18, 5).Value = "Beste klant,`nBedankt voor je e-mail. Kun je me iets meer informatie geven over welk specifiek document of welke informatie je wilt laten doorsturen? Dan kan ik je beter helpen.`nMet vriendelijke groet,`n[Naam Bedrijf] E-mailassistent"
$logs.Cells.Item(18, 6).Value = "2025-07-27 19:45:48"
$logs.Cells.Item(18, 7).Value = "Ja"
$logs.Cells.Item(18, 8).Value = "Nee"
$logs.Cells.Item(18, 9).Value = "Ja"
$logs.Cells.Item(18, 10).Value = "Nee"

# Re-fit the row height (writing a multi-line value otherwise stamps an
# explicit customHeight that the other data rows don't carry).
$logs.Rows.Item(18).EntireRow.AutoFit()

# --- Extend conditional formatting ranges from row 17 -> row 18 -------
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))
$logs.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H18"))
$logs.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I18"))
$logs.Range("J2:J17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J18"))

# --- Update Dashboard summary count ------------------------------------
# "Intern verzoek / Actie voor medewerker" count: 3 -> 4
$dashboard.Cells.Item(4, 2).Value = 4
